$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2031189090006017
$ws.Range("C2").Value = 0.00164196779353612
$ws.Range("B3").Value = 0.2444628943021682
$ws.Range("C3").Value = 0.001927562238199013
